$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2 (1FloodControl)
$ws.Range("B2").Value = 0.3613244253276712
$ws.Range("C2").Value = 0.3613244253276712
$ws.Range("D2").Value = 0.3613244253276711

# Row 3 (2Recreation)
$ws.Range("B3").Value = 0.3776126654937615
$ws.Range("C3").Value = 0.3776126654937615
$ws.Range("D3").Value = 0.3776126654937615

# Row 4 (3Hydroelectric Power)
$ws.Range("B4").Value = 0.2610629091785673
$ws.Range("C4").Value = 0.2610629091785673
$ws.Range("D4").Value = 0.2610629091785673

# Row 5 (1Lo)
$ws.Range("E5").Value = 0.537756771336699
$ws.Range("F5").Value = 0.537756771336699
$ws.Range("G5").Value = 0.537756771336699

# Row 6 (2Med)
$ws.Range("E6").Value = 0.2289086716580906
$ws.Range("F6").Value = 0.2289086716580906
$ws.Range("G6").Value = 0.2289086716580905

# Row 7 (3Hi)
$ws.Range("E7").Value = 0.2333345570052105
$ws.Range("F7").Value = 0.2333345570052105
$ws.Range("G7").Value = 0.2333345570052105
